# Weekly update: insert 6 new price records (date 2022-06-02) at the top
# of the "Manzana" data block (rows 940-945), pushing the existing rows
# 940-986 down to 946-992.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows at 940, shifting the existing 940:986 block down to 946:992.
$ws.Rows("940:945").Insert()

# New rows' data: (row, variety, grade, M, N, O, P, S)
$newRows = @(
    @(940, "Fuji royal",   "Primera", 100, 8000, 9000, 8500, 531),
    @(941, "Fuji royal",   "Segunda", 50,  7000, 7000, 7000, 438),
    @(942, "Granny Smith", "Primera", 100, 8000, 9000, 8500, 531),
    @(943, "Granny Smith", "Segunda", 50,  7000, 7000, 7000, 438),
    @(944, "Pink Lady",    "Primera", 100, 8000, 9000, 8500, 531),
    @(945, "Pink Lady",    "Segunda", 50,  7000, 7000, 7000, 438)
)

foreach ($r in $newRows) {
    $row    = $r[0]
    $variety = $r[1]
    $grade   = $r[2]
    $m = $r[3]
    $n = $r[4]
    $o = $r[5]
    $p = $r[6]
    $s = $r[7]

    $ws.Cells.Item($row, 1).Value = 11
    $ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($row, 3).Value = "Bíobío"
    $ws.Cells.Item($row, 4).Value = 44714
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 8
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100104
    $ws.Cells.Item($row, 8).Value = "Frutos de pepita"
    $ws.Cells.Item($row, 9).Value = 100104002
    $ws.Cells.Item($row, 10).Value = "Manzana"
    $ws.Cells.Item($row, 11).Value = $variety
    $ws.Cells.Item($row, 12).Value = $grade
    $ws.Cells.Item($row, 13).Value = $m
    $ws.Cells.Item($row, 14).Value = $n
    $ws.Cells.Item($row, 15).Value = $o
    $ws.Cells.Item($row, 16).Value = $p
    $ws.Cells.Item($row, 17).Value = "`$/caja 16 kilos empedrada"
    $ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 19).Value = $s
    $ws.Cells.Item($row, 20).Value = 16
}

Write-Output "done"
